# Update the "Equilibrated Data" sheet with new counts/error values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equilibrated Data")

$data = @{
    2  = @{ B = 21.23458333333333;  C = 0.3763829895833333  }
    3  = @{ B = 19.97958333333333;  C = 0.3651268854166666  }
    4  = @{ B = 20.09625;           C = 0.3662541562500001  }
    5  = @{ B = 19.50083333333333;  C = 0.3607654166666666  }
    6  = @{ B = 20.17416666666666;  C = 0.3666654791666666  }
    7  = @{ B = 19.83208333333333;  C = 0.3634229270833333  }
    8  = @{ B = 19.49041666666667;  C = 0.3605727083333334  }
    9  = @{ B = 19.39791666666667;  C = 0.3598313541666667  }
    10 = @{ B = 19.58791666666666;  C = 0.3613970625         }
    11 = @{ B = 18.67208333333333;  C = 0.3533691770833333  }
    12 = @{ B = 19.025;             C = 0.356243125          }
    13 = @{ B = 18.58416666666667;  C = 0.3517053541666667  }
    14 = @{ B = 41.04888888888889;  C = 0.5226891851851851  }
    15 = @{ B = 21.83944444444444;  C = 0.3821902777777778  }
    16 = @{ B = 41.61666666666667;  C = 0.5271444444444444  }
}

foreach ($row in $data.Keys) {
    $ws.Cells.Item($row, 2).Value = $data[$row].B
    $ws.Cells.Item($row, 3).Value = $data[$row].C
}

$ws.Range("A1").Select()

$wb.Save()
